$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.2090032154340836
$ws.Range("C2").Value = 0.5401929260450161
$ws.Range("J2").Value = 0.02572347266881029
$ws.Range("P2").Value = 0.157556270096463
$ws.Range("S2").Value = 0.06752411575562701
$ws.Range("C3").Value = 0.01169590643274854
$ws.Range("J3").Value = 0.01754385964912281
$ws.Range("P3").Value = 0.7192982456140351
$ws.Range("S3").Value = 0.2514619883040936
$ws.Range("J4").Value = 0.06521739130434782
$ws.Range("O4").Value = 0.02173913043478261
$ws.Range("P4").Value = 0.6086956521739131
$ws.Range("S4").Value = 0.3043478260869565
$ws.Range("B6").Value = 0.0546448087431694
$ws.Range("D6").Value = 0.00546448087431694
$ws.Range("F6").Value = 0.06557377049180328
$ws.Range("J6").Value = 0.2677595628415301
$ws.Range("O6").Value = 0.0273224043715847
$ws.Range("Q6").Value = 0.185792349726776
$ws.Range("R6").Value = 0.07103825136612021
$ws.Range("S6").Value = 0.3224043715846995
$ws.Range("B7").Value = 0.105
$ws.Range("D7").Value = 0.02
$ws.Range("F7").Value = 0.035
$ws.Range("J7").Value = 0.145
$ws.Range("O7").Value = 0.01
$ws.Range("Q7").Value = 0.205
$ws.Range("R7").Value = 0.1
$ws.Range("S7").Value = 0.38
$ws.Range("B8").Value = 0.1127450980392157
$ws.Range("D8").Value = 0.02450980392156863
$ws.Range("E8").Value = 0.002450980392156863
$ws.Range("F8").Value = 0.04411764705882353
$ws.Range("J8").Value = 0.1274509803921569
$ws.Range("O8").Value = 0.02205882352941177
$ws.Range("Q8").Value = 0.1446078431372549
$ws.Range("R8").Value = 0.1004901960784314
$ws.Range("S8").Value = 0.4215686274509804
$ws.Range("B9").Value = 0.1173184357541899
$ws.Range("D9").Value = 0.01675977653631285
$ws.Range("F9").Value = 0.06145251396648044
$ws.Range("J9").Value = 0.06145251396648044
$ws.Range("O9").Value = 0.01675977653631285
$ws.Range("Q9").Value = 0.2625698324022346
$ws.Range("R9").Value = 0.08379888268156424
$ws.Range("S9").Value = 0.3798882681564246
$ws.Range("B10").Value = 0.1168300653594771
$ws.Range("D10").Value = 0.0261437908496732
$ws.Range("E10").Value = 0.0008169934640522876
$ws.Range("F10").Value = 0.05800653594771242
$ws.Range("J10").Value = 0.1429738562091503
$ws.Range("O10").Value = 0.0130718954248366
$ws.Range("Q10").Value = 0.2189542483660131
$ws.Range("R10").Value = 0.08006535947712418
$ws.Range("S10").Value = 0.3431372549019608
$ws.Range("G11").Value = 0.1631419939577039
$ws.Range("J11").Value = 0.1057401812688822
$ws.Range("K11").Value = 0.2507552870090635
$ws.Range("L11").Value = 0.4531722054380665
$ws.Range("S11").Value = 0.02719033232628399
$ws.Range("G12").Value = 0.7784810126582279
$ws.Range("J12").Value = 0.1645569620253164
$ws.Range("K12").Value = 0.006329113924050633
$ws.Range("L12").Value = 0.03164556962025317
$ws.Range("S12").Value = 0.0189873417721519
$ws.Range("G13").Value = 0.5833333333333334
$ws.Range("J13").Value = 0.2708333333333333
$ws.Range("S13").Value = 0.1458333333333333
$ws.Range("F15").Value = 0.02790697674418605
$ws.Range("H15").Value = 0.1069767441860465
$ws.Range("I15").Value = 0.07441860465116279
$ws.Range("J15").Value = 0.3813953488372093
$ws.Range("K15").Value = 0.05116279069767442
$ws.Range("M15").Value = 0.009302325581395349
$ws.Range("O15").Value = 0.07441860465116279
$ws.Range("S15").Value = 0.2744186046511628
$ws.Range("F16").Value = 0.02040816326530612
$ws.Range("H16").Value = 0.1173469387755102
$ws.Range("I16").Value = 0.1020408163265306
$ws.Range("J16").Value = 0.413265306122449
$ws.Range("K16").Value = 0.1173469387755102
$ws.Range("M16").Value = 0.02551020408163265
$ws.Range("O16").Value = 0.08163265306122448
$ws.Range("S16").Value = 0.1224489795918367
$ws.Range("F17").Value = 0.01339285714285714
$ws.Range("H17").Value = 0.1674107142857143
$ws.Range("I17").Value = 0.1205357142857143
$ws.Range("J17").Value = 0.4129464285714285
$ws.Range("K17").Value = 0.1138392857142857
$ws.Range("M17").Value = 0.02008928571428572
$ws.Range("O17").Value = 0.06696428571428571
$ws.Range("S17").Value = 0.08482142857142858
$ws.Range("F18").Value = 0.03804347826086957
$ws.Range("H18").Value = 0.1032608695652174
$ws.Range("I18").Value = 0.08152173913043478
$ws.Range("J18").Value = 0.4456521739130435
$ws.Range("K18").Value = 0.09782608695652174
$ws.Range("M18").Value = 0.03260869565217391
$ws.Range("O18").Value = 0.05434782608695652
$ws.Range("S18").Value = 0.1467391304347826
$ws.Range("F19").Value = 0.01843722563652327
$ws.Range("H19").Value = 0.2352941176470588
$ws.Range("I19").Value = 0.06584723441615452
$ws.Range("J19").Value = 0.3590869183494294
$ws.Range("K19").Value = 0.1255487269534679
$ws.Range("M19").Value = 0.02458296751536436
$ws.Range("N19").Value = 0.000877963125548727
$ws.Range("O19").Value = 0.06760316066725197
$ws.Range("S19").Value = 0.1027216856892011
